$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 288.6154
$ws.Range("I6").Value = 187.66667
$ws.Range("K6").Value = 563.00001
$ws.Range("M6").Value = -451.00001

$ws.Range("H41").Value = 1106.8334
$ws.Range("I41").Value = 860.3
$ws.Range("J41").Value = 2339.5
$ws.Range("K41").Value = 860.3
$ws.Range("L41").Value = 2339.5
$ws.Range("M41").Value = -420.3
$ws.Range("N41").Value = -3219.5

$ws.Range("H63").Value = 70271
$ws.Range("J63").Value = 70271
$ws.Range("L63").Value = 70271
$ws.Range("N63").Value = -71519

$ws.Range("H66").Value = 70271
$ws.Range("J66").Value = 70271
$ws.Range("L66").Value = 210813
$ws.Range("N66").Value = -217053

$ws.Range("H98").Value = 1095
$ws.Range("I98").Value = 1095
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 1095
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 403
$ws.Range("N98").ClearContents()

$ws.Range("H107").Value = 9213.333000000001
$ws.Range("J107").Value = 8944
$ws.Range("L107").Value = 8944
$ws.Range("N107").Value = -12784

$ws.Range("H122").Value = 1095
$ws.Range("I122").Value = 1095
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3285
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -835
$ws.Range("N122").ClearContents()

$ws.Range("H125").Value = 3005.75
$ws.Range("I125").Value = 1336
$ws.Range("J125").Value = 3339.7
$ws.Range("K125").Value = 12024
$ws.Range("L125").Value = 30057.3
$ws.Range("M125").Value = -9564
$ws.Range("N125").Value = -34977.3

$ws.Range("H138").Value = 1195027.9
$ws.Range("I138").Value = 1670.8572
$ws.Range("J138").Value = 1493367.1
$ws.Range("K138").Value = 5012.571599999999
$ws.Range("L138").Value = 4480101.300000001
$ws.Range("M138").Value = 127.4284000000007
$ws.Range("N138").Value = -4490381.300000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2535.241
$ws.Range("I32").Value = 1762.3649
$ws.Range("K32").Value = 1762.3649
$ws.Range("M32").Value = -1475.3649

$ws.Range("H33").Value = 5000
$ws.Range("I33").Value = 5000
$ws.Range("K33").Value = 5000
$ws.Range("M33").Value = -4671

$ws.Range("H61").Value = 4438.8096
$ws.Range("I61").Value = 2863.923
$ws.Range("J61").Value = 6998
$ws.Range("K61").Value = 2863.923
$ws.Range("L61").Value = 6998
$ws.Range("M61").Value = -2651.923
$ws.Range("N61").Value = -7422

$ws.Range("H63").Value = 5478.8184
$ws.Range("J63").Value = 7491.25
$ws.Range("L63").Value = 7491.25
$ws.Range("N63").Value = -8863.25

$ws.Range("H66").Value = 5478.8184
$ws.Range("J66").Value = 7491.25
$ws.Range("L66").Value = 37456.25
$ws.Range("N66").Value = -44320.25

$ws.Range("H74").Value = 723.6667
$ws.Range("I74").Value = 423.2857
$ws.Range("K74").Value = 423.2857
$ws.Range("M74").Value = 450.7143

$ws.Range("H77").Value = 723.6667
$ws.Range("I77").Value = 423.2857
$ws.Range("K77").Value = 2116.4285
$ws.Range("M77").Value = 2251.5715

$ws.Range("H101").Value = 53995.5
$ws.Range("J101").Value = 53995.5
$ws.Range("L101").Value = 53995.5
$ws.Range("N101").Value = -60485.5

$ws.Range("H110").Value = 1588.909
$ws.Range("J110").Value = 1628.8334
$ws.Range("L110").Value = 1628.8334
$ws.Range("N110").Value = -5718.8334

$ws.Range("H122").Value = 1999
$ws.Range("I122").Value = 1999
$ws.Range("K122").Value = 5997
$ws.Range("M122").Value = -3547

$ws.Range("H136").Value = 4438.8096
$ws.Range("I136").Value = 2863.923
$ws.Range("J136").Value = 6998
$ws.Range("K136").Value = 8591.769
$ws.Range("L136").Value = 20994
$ws.Range("M136").Value = -6041.769
$ws.Range("N136").Value = -26094

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 4600.3335
$ws.Range("I99").Value = 4650
$ws.Range("J99").Value = 4501
$ws.Range("K99").Value = 4650
$ws.Range("L99").Value = 4501
$ws.Range("M99").Value = -3152
$ws.Range("N99").Value = -7497

$ws.Range("H107").Value = 3461.3572
$ws.Range("I107").Value = 3233.7273
$ws.Range("J107").Value = 4296
$ws.Range("K107").Value = 3233.7273
$ws.Range("L107").Value = 4296
$ws.Range("M107").Value = -1313.7273
$ws.Range("N107").Value = -8136

$ws.Range("H132").Value = 67500
$ws.Range("J132").Value = 67500
$ws.Range("L132").Value = 67500
$ws.Range("N132").Value = -77620

$ws.Range("H134").Value = 4410.6665
$ws.Range("I134").Value = 3924.75
$ws.Range("K134").Value = 11774.25
$ws.Range("M134").Value = -9239.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 5470.75
$ws.Range("I25").Value = 5470.75
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 5470.75
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -5296.75
$ws.Range("N25").ClearContents()

$ws.Range("H59").Value = 43198.1
$ws.Range("I59").Value = 23495
$ws.Range("J59").Value = 48123.875
$ws.Range("K59").Value = 23495
$ws.Range("L59").Value = 48123.875
$ws.Range("M59").Value = -22350
$ws.Range("N59").Value = -50413.875

$ws.Range("H122").Value = 3881.3076
$ws.Range("I122").Value = 2819.625
$ws.Range("K122").Value = 8458.875
$ws.Range("M122").Value = -6008.875

$ws.Range("H134").Value = 6432.7144
$ws.Range("I134").Value = 6117.5557
$ws.Range("K134").Value = 18352.6671
$ws.Range("M134").Value = -15817.6671

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 666.8570999999999
$ws.Range("I14").Value = 666.8570999999999
$ws.Range("K14").Value = 2000.5713
$ws.Range("M14").Value = -1827.5713

$ws.Range("H18").Value = 226
$ws.Range("I18").Value = 226
$ws.Range("K18").Value = 678
$ws.Range("M18").Value = -509

$ws.Range("H21").Value = 256.83334
$ws.Range("I21").Value = 234.2
$ws.Range("K21").Value = 702.5999999999999
$ws.Range("M21").Value = -529.5999999999999

$ws.Range("H68").Value = 2750.3157
$ws.Range("J68").Value = 2886.6572
$ws.Range("L68").Value = 8659.971600000001
$ws.Range("N68").Value = -10281.9716

$ws.Range("H71").Value = 2750.3157
$ws.Range("J71").Value = 2886.6572
$ws.Range("L71").Value = 25979.9148
$ws.Range("N71").Value = -34091.9148

$ws.Range("H98").Value = 294.85715
$ws.Range("I98").Value = 93.333336
$ws.Range("J98").Value = 446
$ws.Range("K98").Value = 280.000008
$ws.Range("L98").Value = 1338
$ws.Range("M98").Value = 1217.999992
$ws.Range("N98").Value = -4334

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()

$ws.Range("H41").Value = 4447.25
$ws.Range("I41").Value = 4447.25
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 4447.25
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -4092.25
$ws.Range("N41").ClearContents()

$ws.Range("H97").Value = 723.9091
$ws.Range("I97").Value = 709.5714
$ws.Range("K97").Value = 709.5714
$ws.Range("M97").Value = -213.5714

$ws.Range("H102").Value = 4014.4666
$ws.Range("I102").Value = 3657.182
$ws.Range("K102").Value = 3657.182
$ws.Range("M102").Value = -2035.182

$ws.Range("H122").Value = 1619.6666
$ws.Range("I122").Value = 1843
$ws.Range("J122").Value = 1396.3334
$ws.Range("K122").Value = 5529
$ws.Range("L122").Value = 4189.0002
$ws.Range("M122").Value = -3079
$ws.Range("N122").Value = -9089.0002

$ws.Range("H132").Value = 8289.833000000001
$ws.Range("I132").Value = 5000
$ws.Range("J132").Value = 8947.799999999999
$ws.Range("K132").Value = 15000
$ws.Range("L132").Value = 26843.4
$ws.Range("M132").Value = -12470
$ws.Range("N132").Value = -31903.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2115.182
$ws.Range("I22").Value = 1718.8572
$ws.Range("J22").Value = 2808.75
$ws.Range("K22").Value = 1718.8572
$ws.Range("L22").Value = 2808.75
$ws.Range("M22").Value = -1423.8572
$ws.Range("N22").Value = -3398.75

$ws.Range("H27").Value = 2115.182
$ws.Range("I27").Value = 1718.8572
$ws.Range("J27").Value = 2808.75
$ws.Range("K27").Value = 1718.8572
$ws.Range("L27").Value = 2808.75
$ws.Range("M27").Value = -1611.8572
$ws.Range("N27").Value = -3022.75

$ws.Range("H46").Value = 2424.45
$ws.Range("I46").Value = 1453.125
$ws.Range("K46").Value = 1453.125
$ws.Range("M46").Value = -1265.125

$ws.Range("H122").Value = 4128.5264
$ws.Range("I122").Value = 2922.1
$ws.Range("J122").Value = 5469
$ws.Range("K122").Value = 8766.299999999999
$ws.Range("L122").Value = 16407
$ws.Range("M122").Value = -6316.299999999999
$ws.Range("N122").Value = -21307

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()

$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()

$ws.Range("H95").Value = 35000
$ws.Range("J95").Value = 35000
$ws.Range("L95").Value = 35000
$ws.Range("N95").Value = -40492

$ws.Range("H122").Value = 5666.706
$ws.Range("I122").Value = 5822.3667
$ws.Range("J122").Value = 4499.25
$ws.Range("K122").Value = 17467.1001
$ws.Range("L122").Value = 13497.75
$ws.Range("M122").Value = -15017.1001
$ws.Range("N122").Value = -18397.75

$ws.Range("H131").Value = 67500
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 67500
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 67500
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = -77580

$ws.Range("H136").Value = 3845.147
$ws.Range("I136").Value = 3901.3076
$ws.Range("J136").Value = 3662.625
$ws.Range("K136").Value = 11703.9228
$ws.Range("L136").Value = 10987.875
$ws.Range("M136").Value = -9153.9228
$ws.Range("N136").Value = -16087.875
